# Applies a new threshold / model-exclusion pass over the EDCR Results
# sheet: recomputed precision/recall/F1 (cols B:D) for every epsilon row,
# plus refreshed NSC/PSC/NRC/PRC counts (cols E:H). Rows 2-66 keep their
# NSC/NRC (cols E/G) at 0 (no excluded model rows yet); rows 67-100 pick
# up the excluded-model counts once the new threshold kicks in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$precision = 0.2385321100917431
$recall    = 0.4262295081967213
$f1        = 0.3058823529411765

# Rows 2-66: PSC/PRC increase, NSC/NRC remain 0.
for ($r = 2; $r -le 66; $r++) {
    $ws.Cells.Item($r, 2).Value = $precision
    $ws.Cells.Item($r, 3).Value = $recall
    $ws.Cells.Item($r, 4).Value = $f1
    $ws.Cells.Item($r, 6).Value = 81
    $ws.Cells.Item($r, 8).Value = 12
}

# Rows 67-100: same precision/recall/F1 and PRC, but NSC/PSC/NRC also move
# because the excluded-model threshold now applies.
for ($r = 67; $r -le 100; $r++) {
    $ws.Cells.Item($r, 2).Value = $precision
    $ws.Cells.Item($r, 3).Value = $recall
    $ws.Cells.Item($r, 4).Value = $f1
    $ws.Cells.Item($r, 5).Value = 28
    $ws.Cells.Item($r, 6).Value = 109
    $ws.Cells.Item($r, 7).Value = 34
    $ws.Cells.Item($r, 8).Value = 12
}
